# Uruguay Primera División — base update, 2024-02-08 00:22
#
# Two pairs of match rows had their data (id, teams, score, result, and all
# odds columns) swapped between rows during re-import: 233<->234 and
# 236<->239. Column A (the running index) is positional and stays put;
# Div / Div Original Name / Date (C:E) are identical within each pair, so
# only columns B and F:AC need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Cell {
    param($ws, [int]$row1, [int]$row2, [int]$col)

    $cell1 = $ws.Cells.Item($row1, $col)
    $cell2 = $ws.Cells.Item($row2, $col)

    $value1 = $cell1.Value2
    $value2 = $cell2.Value2

    $cell1.Value2 = $value2
    $cell2.Value2 = $value1
}

function Swap-MatchRows {
    param($ws, [int]$row1, [int]$row2)

    # Column B = id
    Swap-Cell $ws $row1 $row2 2

    # Columns F:AC = HomeTeam .. PL_AhUnder
    for ($col = 6; $col -le 29; $col++) {
        Swap-Cell $ws $row1 $row2 $col
    }
}

Swap-MatchRows $ws 233 234
Swap-MatchRows $ws 236 239
